# Implemented force param in queue.push_request()
#
# This script renames the "0TestDistributor" / "1TestDistributor" sheets (and the
# matching header labels on the "success"/"fail"/"network" sheets) to
# "0OriginalSequentialForwarding" / "1OriginalSequentialForwarding", refreshes the
# "network" time-series values (and the matching "network" column on the two
# per-distributor sheets) with the new measurements, and appends two new rows of
# data (t=210 and t=215) to every data sheet.

$wb = $excel.ActiveWorkbook

$oldName0 = "0TestDistributor"
$oldName1 = "1TestDistributor"
$newName0 = "0OriginalSequentialForwarding"
$newName1 = "1OriginalSequentialForwarding"

# --- Rename the two worksheet tabs -----------------------------------------
$wsDist0 = $wb.Worksheets.Item($oldName0)
$wsDist0.Name = $newName0

$wsDist1 = $wb.Worksheets.Item($oldName1)
$wsDist1.Name = $newName1

# --- Updated "network" measurements (shared by network/B,C and the two ------
# --- per-distributor sheets' "network" column, column D) -------------------
$networkValues = @{
    2 = 16;   3 = 30;   4 = 46;   5 = 62;   6 = 78;   7 = 94;   8 = 111;  9 = 130;
    10 = 149; 11 = 168; 12 = 191; 13 = 210; 14 = 228; 15 = 246; 16 = 260; 17 = 280;
    18 = 303; 19 = 328; 20 = 354; 21 = 378; 22 = 402; 23 = 423; 24 = 446; 25 = 466;
    26 = 486; 27 = 509; 28 = 533; 29 = 556; 30 = 575; 31 = 594; 32 = 612; 33 = 627;
    34 = 645; 35 = 659; 36 = 668; 37 = 683; 38 = 702; 39 = 722; 40 = 747; 41 = 768;
    42 = 780; 43 = 789; 44 = 792
}

# New trailing time points to append to every sheet.
$newTimes = @{ 43 = 210; 44 = 215 }

# --- success sheet: header rename + two blank (zero) new rows --------------
$wsSuccess = $wb.Worksheets.Item("success")
$wsSuccess.Cells.Item(1, 2).Value = $newName0
$wsSuccess.Cells.Item(1, 3).Value = $newName1
foreach ($r in $newTimes.Keys) {
    $wsSuccess.Cells.Item($r, 1).Value = $newTimes[$r]
    $wsSuccess.Cells.Item($r, 2).Value = 0
    $wsSuccess.Cells.Item($r, 3).Value = 0
}

# --- fail sheet: header rename + two blank (zero) new rows -----------------
$wsFail = $wb.Worksheets.Item("fail")
$wsFail.Cells.Item(1, 2).Value = $newName0
$wsFail.Cells.Item(1, 3).Value = $newName1
foreach ($r in $newTimes.Keys) {
    $wsFail.Cells.Item($r, 1).Value = $newTimes[$r]
    $wsFail.Cells.Item($r, 2).Value = 0
    $wsFail.Cells.Item($r, 3).Value = 0
}

# --- network sheet: header rename + refreshed values + two new rows --------
$wsNetwork = $wb.Worksheets.Item("network")
$wsNetwork.Cells.Item(1, 2).Value = $newName0
$wsNetwork.Cells.Item(1, 3).Value = $newName1
foreach ($r in $networkValues.Keys) {
    $wsNetwork.Cells.Item($r, 2).Value = $networkValues[$r]
    $wsNetwork.Cells.Item($r, 3).Value = $networkValues[$r]
}
foreach ($r in $newTimes.Keys) {
    $wsNetwork.Cells.Item($r, 1).Value = $newTimes[$r]
}

# --- per-distributor sheets: refreshed "network" column + two new rows -----
foreach ($ws in @($wsDist0, $wsDist1)) {
    foreach ($r in $networkValues.Keys) {
        $ws.Cells.Item($r, 4).Value = $networkValues[$r]
    }
    foreach ($r in $newTimes.Keys) {
        $ws.Cells.Item($r, 1).Value = $newTimes[$r]
        $ws.Cells.Item($r, 2).Value = 0
        $ws.Cells.Item($r, 3).Value = 0
    }
}
